$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Collapse the "isn't" gramStart/gramEnd proofErr split back into one run.
# ---------------------------------------------------------------------------
$apos = [char]0x2019
$null = $d.Content.Find.Execute(
    "And this isn" + $apos + "t what the entire team signed up for.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "And this isn" + $apos + "t what the entire team signed up for.", 2)

# ---------------------------------------------------------------------------
# 2) Collapse the "as long as" gramStart/gramEnd proofErr split into one run.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Talk as long as it takes for everyone to contribute and sign the document.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Talk as long as it takes for everyone to contribute and sign the document.", 2)

# ---------------------------------------------------------------------------
# 3) Replace the trailing empty (numbered) paragraph with two blank
#    (un-numbered) paragraphs followed by a new "Test test ..." paragraph
#    that carries the spell-checker proofErr markup, a lastRenderedPageBreak
#    and the _GoBack bookmark, exactly as produced by a real Word save.
# ---------------------------------------------------------------------------
$last = $d.Paragraphs.Last
$last.Range.ListFormat.RemoveNumbers()
$last.Style = "Normal"

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$testXml = "<w:p $wNs>" +
    '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Test </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>test</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>test</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>test</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>test</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>test</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>test</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>test</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>test</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

# Insert the two plain blank paragraphs plus the "Test ..." paragraph in a
# single InsertXML call: it replaces the (now un-numbered, empty) last
# paragraph's own mark with this run of three paragraphs, so the net effect
# is "one trailing empty list item" -> "two blank paragraphs + Test para".
$replacementXml = "<w:p $wNs></w:p><w:p $wNs></w:p>" + $testXml
$null = $last.Range.InsertXML($replacementXml)
